$wb = $excel.ActiveWorkbook

foreach ($name in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1581
    $ws.Range("F4").Value = 101
    $ws.Range("F11").Value = 3796
    $ws.Range("F23").Value = 2600
}
